$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix copy/paste error: header in D1 should read "Cost", not "Amount"
$ws.Range("D1").Value = "Cost"

# Update the active selection to D1 (matches the saved selection state)
$ws.Range("D1").Select()
